$d = $word.ActiveDocument

# --- Step 1: remove the first paragraph entirely -----------------------
# "The Eleventh Virgin  Part II, Chapter I  (second part) ======...="
# This whole paragraph (including its paragraph mark) disappears, so the
# former second paragraph ("By Dorothy Day") becomes paragraph 1.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Delete()

# --- Step 2: turn the old "By Dorothy Day" (bold) paragraph into a ------
# plain-text "% Dorothy Day" paragraph (no bold run formatting at all).
$byLinePara = $d.Paragraphs.Item(1)
$byLineRange = $byLinePara.Range
# Exclude the trailing paragraph mark so only the run text is removed.
$byLineRange.MoveEnd(1, -1) | Out-Null
$byLineRange.Delete()

# Insert fresh, unformatted text at the (now empty) start of the paragraph.
$insertPoint = $d.Range($byLinePara.Range.Start, $byLinePara.Range.Start)
$insertPoint.InsertBefore("% Dorothy Day")
